# Mark the "done" column (F) as TRUE for the rows whose API operations
# were implemented in this commit (api_operations_backlog.xlsx):
#   F2:F9   - api_browser_concept_parents / api_concept_descendants / ...
#   F12:F16 - additional operations completed
#   F33     - api_descriptions_semantic_tags
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$doneRows = @(2, 3, 4, 5, 6, 7, 8, 9, 12, 13, 14, 15, 16, 33)
foreach ($r in $doneRows) {
    $ws.Range("F$r").Formula = "=TRUE()"
}

# Leave the selection where the author ended up editing.
$ws.Range("F33").Select()
